# "Generate Report for Archive"
# - Status text moves from "Ready for handoff" to "In Translation" on every
#   sheet that surfaces it (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - The Status column(s) get narrower to fit the new (shorter) text, same as
#   a report-regeneration pass would do after the content changed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status values (shared-string text "Ready for handoff" -> "In Translation").
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Re-size the Status columns to their new (narrower) fitted width.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
